$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet keeps its 5 rows (header + 4 metadata rows), but a new
# "slug" metadata row is introduced right under the header, pushing
# the previous row 2 (concept/URI mapping), row 3 (medida/dim) and
# row 4 (datatype) content down one row each. The former row 5 - which
# only carried a single "mapping-ano.xlsx" value in column J - is
# replaced entirely by the (now fully populated) datatype row.
#
# Concretely, this is done as plain in-place value overwrites, working
# from the bottom up so no data is clobbered before it is copied down.

# Row 5 only had column J populated before, so columns A-I need the
# same cell formatting (style) that the rest of the table uses before
# we fill them in - copy it down from row 4.
$ws.Range("A4:I4").Copy()
$ws.Range("A5:I5").PasteSpecial(-4122)

# Row 5 (was row 4): datatype row, now populated across every column.
$ws.Range("A5").Value = "null"
$ws.Range("B5").Value = "null"
$ws.Range("C5").Value = "xsd:int"
$ws.Range("D5").Value = "null"
$ws.Range("E5").Value = "URI-Comunidad"
$ws.Range("F5").Value = "xsd:string"
$ws.Range("G5").Value = "xsd:int"
$ws.Range("H5").Value = "xsd:string"
$ws.Range("I5").Value = "URI-comarca"
$ws.Range("J5").Value = "xsd:date"

# Row 4 (was row 3): medida/dim row.
$ws.Range("A4").Value = "null"
$ws.Range("B4").Value = "null"
$ws.Range("C4").Value = "medida"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "dim"
$ws.Range("F4").Value = "medida"
$ws.Range("G4").Value = "medida"
$ws.Range("H4").Value = "medida"
$ws.Range("I4").Value = "dim"
$ws.Range("J4").Value = "dim"

# Row 3 (was row 2): concept / URI mapping row.
$ws.Range("A3").Value = "null"
$ws.Range("B3").Value = "null"
$ws.Range("C3").Value = "iaest-measure:orden"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "iaest-measure:sector-descripcion"
$ws.Range("G3").Value = "iaest-measure:n-accidentes"
$ws.Range("H3").Value = "iaest-measure:temporalidad"
$ws.Range("I3").Value = "sdmx-dimension:refArea"
$ws.Range("J3").Value = "sdmx-dimension:refPeriod"

# Row 2: brand-new machine-readable "slug" identifiers for each column,
# enabling two columns to be related to build SKOS hierarchies.
$ws.Range("A2").Value = "sector-codigo"
$ws.Range("B2").Value = "comarca-codigo"
$ws.Range("C2").Value = "orden"
$ws.Range("D2").Value = "ccaa-codigo"
$ws.Range("E2").Value = "ccaa-nombre"
$ws.Range("F2").Value = "sector-descripcion"
$ws.Range("G2").Value = "n-accidentes"
$ws.Range("H2").Value = "temporalidad"
$ws.Range("I2").Value = "comarca-nombre"
$ws.Range("J2").Value = "ano"
